# BurnDownChart & LogBook update
# - Fills in LogBook rows 38-41 on Sheet1 (Story/Task text + daily effort numbers)
# - Extends the "Ideal" burn-down SUM range in E45 from E6:E37 to E6:E41
# - Updates the sheet view selection to match the new working area

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 38 : Story "Re-Checking, Re-Correct, Re-Designing" /
#              Task "meringkas slide 5 aplikasi media player yang kami sarankan "
$ws.Range("C38").Value = "Re-Checking, Re-Correct, Re-Designing"
$ws.Range("D38").Value = "meringkas slide 5 aplikasi media player yang kami sarankan "
for ($col = 5; $col -le 14; $col++) {
    $ws.Cells.Item(38, $col).Value = 1
}

# --- Row 39 : Story "Tambahan dari Beberapa Aplikasi di android" /
#              Task "Fitur Aplikasi Video Editing VidTrim di Android"
$ws.Range("C39").Value = "Tambahan dari Beberapa Aplikasi di android"
$ws.Range("D39").Value = "Fitur Aplikasi Video Editing VidTrim di Android"
for ($col = 5; $col -le 14; $col++) {
    $ws.Cells.Item(39, $col).Value = 2
}

# --- Row 40 : Story "Tambahan dari Beberapa Aplikasi di android" /
#              Task "Aplikasi penyimpanan data online."
$ws.Range("C40").Value = "Tambahan dari Beberapa Aplikasi di android"
$ws.Range("D40").Value = "Aplikasi penyimpanan data online."
for ($col = 5; $col -le 14; $col++) {
    $ws.Cells.Item(40, $col).Value = 2
}

# --- Row 41 : Story "Re-Checking, Re-Correct, Re-Designing" /
#              Task "Ringkas/Hapus slide 5 Aplikasi Games yang kami sarankan"
$ws.Range("C41").Value = "Re-Checking, Re-Correct, Re-Designing"
$ws.Range("D41").Value = "Ringkas/Hapus slide 5 Aplikasi Games yang kami sarankan"
for ($col = 5; $col -le 14; $col++) {
    $ws.Cells.Item(41, $col).Value = 2
}

# --- Extend the ideal burn-down total to include the newly filled rows
$ws.Range("E45").Formula = "=SUM(E6:E41)"

# --- Update the active view: scroll position + selection
$ws.Activate() | Out-Null
$ws.Range("E41:N41").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 20
